# Remove the header row (number/title/plan) from the Vendor Print Index sheet.
# The upload format no longer writes a header row, so every data row shifts
# up by one and the used range shrinks from A1:C6 to A1:C5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Delete()

# Match the saved selection recorded after the edit.
$ws.Range("C11").Select()
